$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F (ImageLink1 and everything after shifts right by one)
$ws.Columns("F").Insert()

# The newly inserted column inherits the width of the column to its left (Note)
$ws.Columns("F").ColumnWidth = $ws.Columns("E").ColumnWidth

# Header for the new column
$ws.Range("F1").Value = "PdfUpload"

# New column values: whether a PDF upload should be offered for this question
# ("Yes" for subjective questions that only require a typed/handwritten
# definition upload, "No" otherwise)
$ws.Range("F2").Value = "Yes"
$ws.Range("F3").Value = "No"
$ws.Range("F4").Value = "No"
$ws.Range("F5").Value = "No"
$ws.Range("F6").Value = "Yes"
$ws.Range("F7").Value = "No"
$ws.Range("F8").Value = "No"
$ws.Range("F9").Value = "No"

# Match the saved selection state
$ws.Range("F9").Select()
